# Apply "update hoan tra + detail hoa don" changes to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4: 0 -> 1000000000
$ws.Range("F4").Value = 1000000000

# E5: 0 -> 9999
$ws.Range("E5").Value = 9999

# sheetView selection moves from K6 to K7
$ws.Range("K7").Select()

# Column F gets its own explicit width of 11 characters (previously grouped
# with B:XFD at the default 8.88671875). ColumnWidth is in characters; this
# host stores col width as ColumnWidth + 5/6, so 61/6 (~10.1666667) round-trips
# to a stored width of exactly 11, matching the target column F width.
$ws.Range("F:F").ColumnWidth = 10.1666666666667
